# Slide 16, shape "Rectangle 3" (body placeholder) - fix species names / wording
# per the commit "Updated dates and HW keys".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(16)
$sh = $s.Shapes.Item(4)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 2: "bay ancovy (Anchoa mitchilli)," -> "bay anchovy (Anchoa mitchilli),"
#     plus split the genus/species into separate runs.
$para2 = $tr.Paragraphs(2)

# "ancovy " (7 chars, incl. trailing space) -> "anchovy " ; leaves "(" as its own run
$para2.Characters(5, 7).Text = "anchovy "

# Split the (already-italic) "Anchoa mitchilli" run into three runs, preserving
# italics, by re-assigning each sub-range's text to itself.
$para2.Characters(14, 6).Text = "Anchoa"      # "Anchoa"
$para2.Characters(20, 1).Text = " "           # " "
$para2.Characters(21, 9).Text = "mitchilli"   # "mitchilli"

# --- Paragraph 3: "bluefish (Pomatomus saltatrix)," -> split genus/species into runs
$para3 = $tr.Paragraphs(3)
$para3.Characters(11, 9).Font.Italic = $true   # "Pomatomus"
$para3.Characters(20, 1).Font.Italic = $true   # " "
$para3.Characters(21, 9).Font.Italic = $true   # "saltatrix"

# --- Paragraph 4: "striped bass (Morone saxatilis), and" -> split genus/species into runs
$para4 = $tr.Paragraphs(4)
$para4.Characters(15, 6).Font.Italic = $true   # "Morone"
$para4.Characters(21, 1).Font.Italic = $true   # " "
$para4.Characters(22, 9).Font.Italic = $true   # "saxatilis"

# --- Paragraph 5: "weakfish (Cynoscion regalis)." -> split genus/species into runs
$para5 = $tr.Paragraphs(5)
$para5.Characters(11, 9).Font.Italic = $true   # "Cynoscion"
$para5.Characters(20, 1).Font.Italic = $true   # " "
$para5.Characters(21, 7).Font.Italic = $true   # "regalis"
